$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-02 Thursday" "2025-01-03 Friday"

Replace-Text "96÷6=" "44÷7="
Replace-Text "94÷8=" "72÷2="
Replace-Text "90÷5=" "61÷4="
Replace-Text "91÷8=" "96÷9="
Replace-Text "87÷4=" "24÷6="

Replace-Text "43÷3=" "97÷8="
Replace-Text "73÷8=" "12÷7="
Replace-Text "49÷6=" "32÷8="
Replace-Text "37÷5=" "34÷4="
Replace-Text "25÷3=" "20÷7="

Replace-Text "63÷8=" "81÷9="
Replace-Text "81÷5=" "51÷7="
Replace-Text "69÷7=" "52÷4="
Replace-Text "76÷8=" "38÷3="
Replace-Text "49÷2=" "81÷6="

Replace-Text "61÷7=" "24÷3="
Replace-Text "54÷5=" "87÷3="
Replace-Text "42÷2=" "90÷4="
Replace-Text "68÷2=" "90÷5="
Replace-Text "54÷9=" "36÷8="

Replace-Text "35÷9=" "63÷2="
Replace-Text "52÷3=" "68÷9="
Replace-Text "28÷5=" "54÷3="
Replace-Text "85÷6=" "12÷5="
Replace-Text "52÷2=" "22÷6="
